# Update column G ("K") values on the active sheet to reflect the
# regenerated std/mean-based s_vals calculation (commit: "regen save_data
# to use K instead of Strike#, regen std/mean, calc and write s_vals").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 2
    4  = 5
    5  = 1
    6  = 1
    7  = 0
    8  = 3
    9  = 2
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 2
    15 = 2
    16 = 2
    17 = 1
    18 = 0
    19 = 2
    20 = 0
    21 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
